$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 263396
$ws.Range("J33").Value = 396.33334
$ws.Range("L33").Value = 396.33334
$ws.Range("N33").Value = -854.33334
$ws.Range("H88").Value = 20010798
$ws.Range("I88").Value = 200000000
$ws.Range("J88").Value = 11998.777
$ws.Range("K88").Value = 200000000
$ws.Range("L88").Value = 11998.777
$ws.Range("M88").Value = -199999594
$ws.Range("N88").Value = -12810.777
$ws.Range("H91").Value = 20010798
$ws.Range("I91").Value = 200000000
$ws.Range("J91").Value = 11998.777
$ws.Range("K91").Value = 200000000
$ws.Range("L91").Value = 11998.777
$ws.Range("M91").Value = -199998596
$ws.Range("N91").Value = -14806.777
$ws.Range("H107").Value = 720.7778
$ws.Range("I107").Value = 685.125
$ws.Range("J107").Value = 1006
$ws.Range("K107").Value = 685.125
$ws.Range("L107").Value = 1006
$ws.Range("M107").Value = 1234.875
$ws.Range("N107").Value = -4846
$ws.Range("H116").Value = 39590240
$ws.Range("J116").Value = 25005028
$ws.Range("L116").Value = 25005028
$ws.Range("N116").Value = -25011912
$ws.Range("H129").Value = 1407.6086
$ws.Range("I129").Value = 802.5
$ws.Range("J129").Value = 1730.3334
$ws.Range("K129").Value = 2407.5
$ws.Range("L129").Value = 5191.0002
$ws.Range("M129").Value = 2592.5
$ws.Range("N129").Value = -15191.0002
$ws.Range("H137").Value = 8337271.5
$ws.Range("J137").Value = 11910042
$ws.Range("L137").Value = 35730126
$ws.Range("N137").Value = -35735226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6113.3477
$ws.Range("I32").Value = 3372.4324
$ws.Range("K32").Value = 3372.4324
$ws.Range("M32").Value = -3085.4324
$ws.Range("H102").Value = 372035.34
$ws.Range("I102").Value = 762253.75
$ws.Range("K102").Value = 762253.75
$ws.Range("M102").Value = -760631.75
$ws.Range("H132").Value = 4838.7104
$ws.Range("I132").Value = 2028.2084
$ws.Range("J132").Value = 9656.714
$ws.Range("K132").Value = 6084.6252
$ws.Range("L132").Value = 28970.142
$ws.Range("M132").Value = -3554.6252
$ws.Range("N132").Value = -34030.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1681.7
$ws.Range("I107").Value = 1659.4706
$ws.Range("J107").Value = 1807.6666
$ws.Range("K107").Value = 1659.4706
$ws.Range("L107").Value = 1807.6666
$ws.Range("M107").Value = 260.5293999999999
$ws.Range("N107").Value = -5647.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5800.936
$ws.Range("I31").Value = 1894.5238
$ws.Range("K31").Value = 1894.5238
$ws.Range("M31").Value = -1599.5238
$ws.Range("H34").Value = 5800.936
$ws.Range("I34").Value = 1894.5238
$ws.Range("K34").Value = 1894.5238
$ws.Range("M34").Value = -1692.5238
$ws.Range("H94").Value = 2801.45
$ws.Range("I94").Value = 2745.5557
$ws.Range("J94").Value = 2847.182
$ws.Range("K94").Value = 2745.5557
$ws.Range("L94").Value = 2847.182
$ws.Range("M94").Value = -2294.5557
$ws.Range("N94").Value = -3749.182
$ws.Range("I99").Value = 2647.5715
$ws.Range("J99").Value = 10571.429
$ws.Range("K99").Value = 2647.5715
$ws.Range("L99").Value = 10571.429
$ws.Range("M99").Value = -1149.5715
$ws.Range("N99").Value = -13567.429
$ws.Range("H105").Value = 1370.5
$ws.Range("J105").Value = 1495.5
$ws.Range("L105").Value = 1495.5
$ws.Range("N105").Value = -4989.5
$ws.Range("H122").Value = 4159.4443
$ws.Range("I122").Value = 2732.1
$ws.Range("K122").Value = 8196.299999999999
$ws.Range("M122").Value = -5746.299999999999
$ws.Range("I126").Value = 2647.5715
$ws.Range("J126").Value = 10571.429
$ws.Range("K126").Value = 7942.7145
$ws.Range("L126").Value = 31714.287
$ws.Range("M126").Value = -5472.7145
$ws.Range("N126").Value = -36654.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 93.14815
$ws.Range("I2").Value = 77.818184
$ws.Range("J2").Value = 103.6875
$ws.Range("K2").Value = 466.909104
$ws.Range("L2").Value = 622.125
$ws.Range("M2").Value = -353.909104
$ws.Range("N2").Value = -848.125
$ws.Range("H129").Value = 1534.5
$ws.Range("J129").Value = 1712.6666
$ws.Range("L129").Value = 5137.9998
$ws.Range("N129").Value = -15137.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 408.5
$ws.Range("I97").Value = 390
$ws.Range("K97").Value = 390
$ws.Range("M97").Value = 106
$ws.Range("H122").Value = 6221.857
$ws.Range("I122").Value = 2785
$ws.Range("J122").Value = 7596.6
$ws.Range("K122").Value = 8355
$ws.Range("L122").Value = 22789.8
$ws.Range("M122").Value = -5905
$ws.Range("N122").Value = -27689.8
$ws.Range("H126").Value = 4342.84
$ws.Range("J126").Value = 5677.6665
$ws.Range("L126").Value = 17032.9995
$ws.Range("N126").Value = -21972.9995
$ws.Range("H132").Value = 5158.1875
$ws.Range("I132").Value = 5437.222
$ws.Range("J132").Value = 4799.4287
$ws.Range("K132").Value = 16311.666
$ws.Range("L132").Value = 14398.2861
$ws.Range("M132").Value = -13781.666
$ws.Range("N132").Value = -19458.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5246.4165
$ws.Range("I7").Value = 1984
$ws.Range("J7").Value = 6333.8887
$ws.Range("K7").Value = 1984
$ws.Range("L7").Value = 6333.8887
$ws.Range("M7").Value = -1872
$ws.Range("N7").Value = -6557.8887
$ws.Range("H22").Value = 1217.1
$ws.Range("I22").Value = 994.6
$ws.Range("K22").Value = 994.6
$ws.Range("M22").Value = -699.6
$ws.Range("H27").Value = 1217.1
$ws.Range("I27").Value = 994.6
$ws.Range("K27").Value = 994.6
$ws.Range("M27").Value = -887.6
$ws.Range("H40").Value = 23810688
$ws.Range("I40").Value = 1251.9231
$ws.Range("K40").Value = 1251.9231
$ws.Range("M40").Value = -1115.9231
$ws.Range("H122").Value = 57150984
$ws.Range("I122").Value = 83337050
$ws.Range("K122").Value = 250011150
$ws.Range("M122").Value = -250008700
$ws.Range("H126").Value = 5246.4165
$ws.Range("I126").Value = 1984
$ws.Range("J126").Value = 6333.8887
$ws.Range("K126").Value = 5952
$ws.Range("L126").Value = 19001.6661
$ws.Range("M126").Value = -3482
$ws.Range("N126").Value = -23941.6661
$ws.Range("H132").Value = 6826.3706
$ws.Range("I132").Value = 4837.2
$ws.Range("J132").Value = 7996.4707
$ws.Range("K132").Value = 14511.6
$ws.Range("L132").Value = 23989.4121
$ws.Range("M132").Value = -11981.6
$ws.Range("N132").Value = -29049.4121
$ws.Range("H133").Value = 92326
$ws.Range("J133").Value = 92326
$ws.Range("L133").Value = 92326
$ws.Range("N133").Value = -97386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3817.8333
$ws.Range("J96").Value = 4541.4
$ws.Range("L96").Value = 4541.4
$ws.Range("N96").Value = -7287.4
$ws.Range("H116").Value = 91141.75
$ws.Range("J116").Value = 91141.75
$ws.Range("L116").Value = 91141.75
$ws.Range("N116").Value = -100319.75
$ws.Range("H122").Value = 5336.222
$ws.Range("I122").Value = 4575.1904
$ws.Range("K122").Value = 13725.5712
$ws.Range("M122").Value = -11275.5712
$ws.Range("H123").Value = 92429
$ws.Range("J123").Value = 92429
$ws.Range("L123").Value = 92429
$ws.Range("N123").Value = -102229
$ws.Range("H125").Value = 65000
$ws.Range("J125").Value = 65000
$ws.Range("L125").Value = 65000
$ws.Range("N125").Value = -74840
$ws.Range("H126").Value = 2078.5557
$ws.Range("I126").Value = 1965
$ws.Range("J126").Value = 2257
$ws.Range("K126").Value = 5895
$ws.Range("L126").Value = 6771
$ws.Range("M126").Value = -3425
$ws.Range("N126").Value = -11711
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
$ws.Range("H132").Value = 55564576
$ws.Range("I132").Value = 9260810
$ws.Range("K132").Value = 27782430
$ws.Range("M132").Value = -27779900
$ws.Range("H133").Value = 60139
$ws.Range("J133").Value = 60139
$ws.Range("L133").Value = 60139
$ws.Range("N133").Value = -70259
$ws.Range("H136").Value = 8532.547
$ws.Range("I136").Value = 3411.7585
$ws.Range("K136").Value = 10235.2755
$ws.Range("M136").Value = -7685.2755
$ws.Range("H137").Value = 90715
$ws.Range("J137").Value = 90715
$ws.Range("L137").Value = 90715
$ws.Range("N137").Value = -100915
$ws.Range("H139").Value = 55555
$ws.Range("J139").Value = 55555
$ws.Range("L139").Value = 55555
$ws.Range("N139").Value = -65835
$ws.Range("H141").Value = 49999
$ws.Range("J141").Value = 49998
$ws.Range("L141").Value = 49998
$ws.Range("N141").Value = -60358
